$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
    3  = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg"
    4  = "gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    5  = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
    6  = "alshimaa.atef@med.asu.edu.egm, mennatulla.medhat@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
    7  = "AbeerRagheb@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"
    8  = "AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
    9  = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    11 = "aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    12 = "yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
    13 = "esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
    15 = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
    17 = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"
    20 = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
    25 = "menna-allah.gamil@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg"
    27 = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    28 = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
    30 = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
